$d = $word.ActiveDocument

# ===========================================================================
# Section 1: bioRxiv preprint -> Ivan Maslov et al., Communications Biology
# ===========================================================================

# --- Paragraph 1: date / venue line --------------------------------------
$null = $d.Content.Find.Execute(
    "Co-Author, Publishing date: 2020-01-01, bioRxiv (recently accepted at Communications Biology)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Co-Author, Publishing date: 2023-03-01, Communications Biology",
    2)

# --- Paragraph 2: authors line (contains literal straight apostrophes that
#     must NOT be smart-quoted, so assign Range.Text directly). The leading
#     "Authors" word (and its spell-check proofErr wrapper) is left alone so
#     only the part starting at the opening quote is swapped out. ----------
$rng = $d.Content
$found = $rng.Find.Execute(
    "'Danai Laskaratou, Guillermo Solís Fernández, Quinten Coucke, Eduard Fron, Susana Rocha, Johan Hofkens, Jelle Hendrix & Hideaki Mizuno'",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)
if ($found) {
    $rng.Text = "‘Ivan Maslov, Oleksandr Volkov, Polina Khorn, Philipp Orekhov, Anastasiia Gusach, Pavel Kuzmichev, Andrey Gerasimov, Aleksandra Luginina, Quinten Coucke, Andrey Bogorodskiy, Valentin Gordeliy, Simon Wanninger, Anders Barth, Alexey Mishin, Johan Hofkens, Vadim Cherezov, Thomas Gensch, Jelle Hendrix, Valentin Borshchevskiy' "
}

# --- Paragraph 3: DOI line (display text + underlying hyperlink target) --
$null = $d.Content.Find.Execute(
    "https://doi.org/10.1101/2020.11.26.40018",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://doi.org/10.1038/s42003-023-04727-z",
    2)

foreach ($h in $d.Hyperlinks) {
    if ($h.Address -eq "https://doi.org/10.1101/2020.11.26.40018") {
        $h.Address = "https://doi.org/10.1038/s42003-023-04727-z"
        $h.Range.Style = "Hyperlink"
    }
}

# ===========================================================================
# Section 2: bioRxiv preprint -> Hongbo Yuan et al., PNAS
# ===========================================================================

# --- Paragraph 1: date / venue line --------------------------------------
$null = $d.Content.Find.Execute(
    "Co-Author, Publishing date: 2022-08-22, bioRxiv (now accepted at PLOS, expected summer ’23)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Co-Author, Publishing date: 2023-03-03, PNAS",
    2)

# --- Paragraph 3: DOI line (display text + underlying hyperlink target),
#     plus a trailing space added after the hyperlink -----------------------
$null = $d.Content.Find.Execute(
    "https://doi.org/10.1101/2022.08.24.505064",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://doi.org/10.1073/pnas.2216934120",
    2)

foreach ($h in $d.Hyperlinks) {
    if ($h.Address -eq "https://doi.org/10.1101/2022.08.24.505064") {
        $h.Address = "https://doi.org/10.1073/pnas.2216934120"
        $h.Range.Style = "Hyperlink"
    }
}

$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "https://doi.org/10.1073/pnas.2216934120",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $rng2.InsertAfter(" ")
}
